$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-51: refresh Price (D) and Volume(1h) (E) snapshot values.
#
# Some of the new Price strings (e.g. "1.00", "0.998") are plain decimals
# that Excels normal auto-type-detection would coerce into numbers (losing
# the exact text, e.g. trailing zeros or switching to scientific notation).
# Force those specific cells to Text format before assigning so the literal
# string is kept, then restore the Normal style so no stray formatting is
# left behind on the cell.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '629.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.96'
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '466.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.708'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.978'
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '157.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000279'
$ws.Range("D51").Style = "Normal"

# Remaining cells (percentages and multi-dot "thousands" prices) already
# round-trip as text, so a direct assignment is enough.
$ws.Range("D2").Value = '69.197.47'
$ws.Range("E2").Value = '  +2.33%  '
$ws.Range("D3").Value = '3.815.31'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("E5").Value = '  +5.26%  '
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("D7").Value = '3.814.44'
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +2.73%  '
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("E12").Value = '  +3.34%  '
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").Value = '4.451.72'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").Value = '3.818.08'
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").Value = '69.151.77'
$ws.Range("E17").Value = '  +2.26%  '
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E21").Value = '  +1.36%  '
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("E26").Value = '  +3.48%  '
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '3.961.61'
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("E31").Value = '  +3.78%  '
$ws.Range("E32").Value = '  +1.39%  '
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  +1.36%  '
$ws.Range("E37").Value = '  +3.72%  '
$ws.Range("E38").Value = '  +8.14%  '
$ws.Range("E39").Value = '  +6.82%  '
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  +3.67%  '
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("E46").Value = '  +5.82%  '
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("E49").Value = '  +3.39%  '
$ws.Range("E50").Value = '  +1.50%  '
$ws.Range("E51").Value = '  +13.23%  '
